# Auto-generated: apply Kraken_Profits price-refresh values per sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H15").Value = 104.07692
$ws.Range("I15").Value = 104.07692
$ws.Range("K15").Value = 312.23076
$ws.Range("M15").Value = -143.23076
$ws.Range("H40").Value = 6327.3794
$ws.Range("I40").Value = 3915.3333
$ws.Range("J40").Value = 6605.6924
$ws.Range("K40").Value = 3915.3333
$ws.Range("L40").Value = 6605.6924
$ws.Range("M40").Value = -3740.3333
$ws.Range("N40").Value = -6955.6924
$ws.Range("H96").Value = 4390
$ws.Range("I96").Value = 4237.5
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 12712.5
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = -11339.5
$ws.Range("N96").Value = -17746
$ws.Range("H132").Value = 3177.1724
$ws.Range("J132").Value = 5733.3335
$ws.Range("L132").Value = 17200.0005
$ws.Range("N132").Value = -22260.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H24").Value = 45427
$ws.Range("J24").Value = 45427
$ws.Range("L24").Value = 45427
$ws.Range("N24").Value = -46175
$ws.Range("H61").Value = 3666.6667
$ws.Range("I61").Value = 3666.6667
$ws.Range("K61").Value = 3666.6667
$ws.Range("M61").Value = -3454.6667
$ws.Range("H74").Value = 4399.6
$ws.Range("I74").Value = 1666
$ws.Range("K74").Value = 1666
$ws.Range("M74").Value = -792
$ws.Range("H77").Value = 4399.6
$ws.Range("I77").Value = 1666
$ws.Range("K77").Value = 8330
$ws.Range("M77").Value = -3962
$ws.Range("H92").Value = 49999.5
$ws.Range("J92").Value = 49999.5
$ws.Range("L92").Value = 49999.5
$ws.Range("N92").Value = -54991.5
$ws.Range("H96").Value = 12000
$ws.Range("J96").Value = 12000
$ws.Range("L96").Value = 12000
$ws.Range("N96").Value = -17492
$ws.Range("H100").Value = 45427
$ws.Range("J100").Value = 45427
$ws.Range("L100").Value = 45427
$ws.Range("N100").Value = -47591
$ws.Range("H101").Value = 24999.666
$ws.Range("J101").Value = 24999.666
$ws.Range("L101").Value = 24999.666
$ws.Range("N101").Value = -31489.666
$ws.Range("H110").Value = 944.44446
$ws.Range("I110").Value = 962.5
$ws.Range("J110").Value = 800
$ws.Range("K110").Value = 962.5
$ws.Range("L110").Value = 800
$ws.Range("M110").Value = 1082.5
$ws.Range("N110").Value = -4890
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 2959.7144
$ws.Range("I132").Value = 2543.6
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7630.799999999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5100.799999999999
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 3666.6667
$ws.Range("I136").Value = 3666.6667
$ws.Range("K136").Value = 11000.0001
$ws.Range("M136").Value = -8450.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H20").Value = 825
$ws.Range("I20").Value = 1200
$ws.Range("J20").Value = 450
$ws.Range("K20").Value = 1200
$ws.Range("L20").Value = 450
$ws.Range("M20").Value = -953
$ws.Range("N20").Value = -944
$ws.Range("H95").Value = 18541
$ws.Range("J95").Value = 18541
$ws.Range("L95").Value = 18541
$ws.Range("N95").Value = -24033
$ws.Range("H100").Value = 8599.6
$ws.Range("J100").Value = 8599.6
$ws.Range("L100").Value = 8599.6
$ws.Range("N100").Value = -10763.6
$ws.Range("H103").Value = 3476
$ws.Range("J103").Value = 3476
$ws.Range("L103").Value = 3476
$ws.Range("N103").Value = -5820
$ws.Range("H134").Value = 5673
$ws.Range("I134").Value = 3552.75
$ws.Range("J134").Value = 8500
$ws.Range("K134").Value = 10658.25
$ws.Range("L134").Value = 25500
$ws.Range("M134").Value = -8123.25
$ws.Range("N134").Value = -30570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3994.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3994.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3994.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -6240.5
$ws.Range("H89").Value = 3994.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3994.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 19972.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -31204.5
$ws.Range("H132").Value = 2178.8
$ws.Range("I132").Value = 2178.8
$ws.Range("K132").Value = 6536.400000000001
$ws.Range("M132").Value = -4006.400000000001
$ws.Range("H134").Value = 1087.75
$ws.Range("I134").Value = 1087.6666
$ws.Range("J134").Value = 1088
$ws.Range("K134").Value = 3262.9998
$ws.Range("L134").Value = 3264
$ws.Range("M134").Value = -727.9998000000001
$ws.Range("N134").Value = -8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 276.7143
$ws.Range("I12").Value = 435.25
$ws.Range("K12").Value = 1305.75
$ws.Range("M12").Value = -1132.75
$ws.Range("H51").Value = 575
$ws.Range("I51").Value = 575
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1725
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -1265
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 1750
$ws.Range("I58").Value = 1500
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 4500
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -4372
$ws.Range("N58").Value = -6256
$ws.Range("H133").Value = 15000
$ws.Range("I133").Value = 15000
$ws.Range("K133").Value = 45000
$ws.Range("M133").Value = -39940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3416.6667
$ws.Range("I80").Value = 2875
$ws.Range("K80").Value = 2875
$ws.Range("M80").Value = -1877
$ws.Range("H83").Value = 3416.6667
$ws.Range("I83").Value = 2875
$ws.Range("K83").Value = 14375
$ws.Range("M83").Value = -9383
$ws.Range("H102").Value = 2547.3845
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1625
$ws.Range("I122").Value = 1625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2425
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 8497.571
$ws.Range("I132").Value = 7996
$ws.Range("J132").Value = 9166.333000000001
$ws.Range("K132").Value = 23988
$ws.Range("L132").Value = 27498.999
$ws.Range("M132").Value = -21458
$ws.Range("N132").Value = -32558.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8428.286
$ws.Range("I7").Value = 7749.75
$ws.Range("K7").Value = 7749.75
$ws.Range("M7").Value = -7637.75
$ws.Range("H46").Value = 4561.7354
$ws.Range("I46").Value = 3147.25
$ws.Range("J46").Value = 4750.3335
$ws.Range("K46").Value = 3147.25
$ws.Range("L46").Value = 4750.3335
$ws.Range("M46").Value = -2959.25
$ws.Range("N46").Value = -5126.3335
$ws.Range("H61").Value = 8000
$ws.Range("I61").Value = 8000
$ws.Range("K61").Value = 8000
$ws.Range("M61").Value = -7798
$ws.Range("H113").Value = 8000
$ws.Range("I113").Value = 8000
$ws.Range("K113").Value = 8000
$ws.Range("M113").Value = -5830
$ws.Range("H126").Value = 8428.286
$ws.Range("I126").Value = 7749.75
$ws.Range("K126").Value = 23249.25
$ws.Range("M126").Value = -20779.25
$ws.Range("H136").Value = 5324
$ws.Range("I136").Value = 4986.25
$ws.Range("K136").Value = 14958.75
$ws.Range("M136").Value = -12408.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 95243.75
$ws.Range("J86").Value = 95243.75
$ws.Range("L86").Value = 95243.75
$ws.Range("N86").Value = -97489.75
$ws.Range("H89").Value = 95243.75
$ws.Range("J89").Value = 95243.75
$ws.Range("L89").Value = 476218.75
$ws.Range("N89").Value = -487450.75
$ws.Range("H107").Value = 3066.3333
$ws.Range("I107").Value = 600
$ws.Range("K107").Value = 1800
$ws.Range("M107").Value = 120
$ws.Range("H126").Value = 2571.6
$ws.Range("I126").Value = 2571.6
$ws.Range("K126").Value = 7714.799999999999
$ws.Range("M126").Value = -5244.799999999999
$ws.Range("H132").Value = 984
$ws.Range("I132").Value = 984
$ws.Range("K132").Value = 2952
$ws.Range("M132").Value = -422
$ws.Range("H136").Value = 4464.154
$ws.Range("I136").Value = 4002.8333
$ws.Range("K136").Value = 12008.4999
$ws.Range("M136").Value = -9458.499899999999
